$d = $word.ActiveDocument

# Locate the signature line "CAIO MARINHO DO REIS".
$find = $d.Content.Find
$found = $find.Execute("CAIO MARINHO DO REIS", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $nameRange = $find.Parent
    $nameEnd = $nameRange.End

    # Re-key the trailing "S" as its own edit (mirrors the author retyping the
    # last letter of the surname), which makes Word split it into its own run
    # while keeping identical run formatting (rFonts eastAsia Arial, bold, szCs 24).
    $lastChar = $d.Range($nameEnd - 1, $nameEnd)
    $lastChar.Font.Bold = 0
    $lastChar.Font.Bold = 1
}
